$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Baltimore as a new row (row 6) under the existing cities table
$ws.Range("A6").Value = "Baltimore"
$ws.Range("B6").Value = 39.2904
$ws.Range("C6").Value = -76.6122

# Update the active selection to match the target workbook state
$ws.Range("K10").Select()
